# Changed motor pins due to not needed motorDirPins.
#
# The unused motorDirPins rows (D44 motor4aPin / D45 motor4bPin) are
# removed, and the motor2/motor3/motor4 a/b pin labels shift up by two
# rows (D36..D41) to fill the gap left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D36").Value = "motor2aPin"
$ws.Range("D37").Value = "motor2bPin"
$ws.Range("D38").Value = "motor3aPin"
$ws.Range("D39").Value = "motor3bPin"
$ws.Range("D40").Value = "motor4aPin"
$ws.Range("D41").Value = "motor4bPin"
$ws.Range("D44").Value = ""
$ws.Range("D45").Value = ""

# Reflect the resulting scroll / selection change in the frozen-pane view.
$ws.Activate()
$ws.Range("D34").Select()
